$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 123
$ws.Range("D3").Value = 124
$ws.Range("D4").Value = 125
$ws.Range("D5").Value = 126

$ws.Range("G16").Select()
